$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 97; existing rows 97..133 shift down to 98..134.
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new weekly price record.
$ws.Cells.Item(97, 1).Value  = 10
$ws.Cells.Item(97, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(97, 3).Value  = "La Araucanía"
$ws.Cells.Item(97, 4).Value  = 45229
$ws.Cells.Item(97, 5).Value  = 9
$ws.Cells.Item(97, 6).Value  = "Fruta"
$ws.Cells.Item(97, 7).Value  = 100108
$ws.Cells.Item(97, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(97, 9).Value  = 100108004
$ws.Cells.Item(97, 10).Value = "Papaya"
$ws.Cells.Item(97, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(97, 12).Value = "Primera"
$ws.Cells.Item(97, 13).Value = 140
$ws.Cells.Item(97, 14).Value = 24000
$ws.Cells.Item(97, 15).Value = 24000
$ws.Cells.Item(97, 16).Value = 24000
$ws.Cells.Item(97, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(97, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(97, 19).Value = 2400
$ws.Cells.Item(97, 20).Value = 10

# Make sure the date cell keeps the same date number format used by the rest of column D.
$ws.Cells.Item(97, 4).NumberFormat = $ws.Cells.Item(98, 4).NumberFormat
